$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old meta_title / meta_description headers (columns M1, N1) and
# replace with the new header set starting at M1.
$ws.Range("M1").Value = "discount"
$ws.Range("N1").Value = "discount_type"
$ws.Range("O1").Value = "part_no"
$ws.Range("P1").Value = "segment_id"
$ws.Range("Q1").Value = "model_id"
$ws.Range("R1").Value = "smart_part_no"
$ws.Range("S1").Value = "ref_part_no"
$ws.Range("T1").Value = "oe_part_no"
$ws.Range("U1").Value = "size"
$ws.Range("V1").Value = "mega_categories"
$ws.Range("W1").Value = "series"

# Apply the same header style (bold) used by the existing headers to the
# newly added header cells.
$ws.Range("A1").Copy()
$ws.Range("O1:W1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 sample data
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = "percentage/amount"
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2

# Set column N width (for discount_type)
$ws.Columns.Item(14).ColumnWidth = 18

# Update selection to match target state
$ws.Range("R2").Select()
